$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new TODO row (row 4) with the same two-column layout as the
# existing rows: task description in column A, status in column B.
$ws.Range("A4").Value = "Try to run PCA or ICA before doing source seperation"
$ws.Range("B4").Value = "waiting"

# Move the selection down to the next empty row, mirroring the original
# sheet's convention of keeping the cursor on the first unused row.
$ws.Range("A5").Select()
